$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.011618333333334
$ws.Range("H2").Value = 3.034855
$ws.Range("I2").Value = 0.5235149663433657
$ws.Range("J2").Value = 0.5235149663433657
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.315861666666667
$ws.Range("N2").Value = 3.947585
$ws.Range("O2").Value = 0.2754050739440597
$ws.Range("P2").Value = 0.2754050739440597
$ws.Range("Q2").Value = 1.331149786130556
$ws.Range("R2").Value = 11.980348075175
$ws.Range("S2").Value = 0.1441786780166166
$ws.Range("T2").Value = 0.1441786780166166
$ws.Range("G3").Value = 1.011618333333334
$ws.Range("H3").Value = 3.034855
$ws.Range("I3").Value = 0.5235149663433657
$ws.Range("J3").Value = 0.5235149663433657
$ws.Range("O3").Value = 0.3040809095127364
$ws.Range("P3").Value = 0.3040809095127364
$ws.Range("Q3").Value = 1.469752288392778
$ws.Range("R3").Value = 13.227770595535
$ws.Range("S3").Value = 0.1591909071092202
$ws.Range("T3").Value = 0.1591909071092202
$ws.Range("G4").Value = 1.011618333333334
$ws.Range("H4").Value = 3.034855
$ws.Range("I4").Value = 0.5235149663433657
$ws.Range("J4").Value = 0.5235149663433657
$ws.Range("M4").Value = 2.009179666666667
$ws.Range("N4").Value = 6.027539
$ws.Range("O4").Value = 0.4205140165432039
$ws.Range("P4").Value = 0.4205140165432039
$ws.Range("Q4").Value = 2.032522985760556
$ws.Range("R4").Value = 18.292706871845
$ws.Range("S4").Value = 0.2201453812175289
$ws.Range("T4").Value = 0.2201453812175289
$ws.Range("I5").Value = 0.2899264353016711
$ws.Range("J5").Value = 0.2899264353016712
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.315861666666667
$ws.Range("N5").Value = 3.947585
$ws.Range("O5").Value = 0.2754050739440597
$ws.Range("P5").Value = 0.2754050739440597
$ws.Range("Q5").Value = 0.7372005332361111
$ws.Range("R5").Value = 6.634804799125
$ws.Range("S5").Value = 0.07984721135259439
$ws.Range("T5").Value = 0.0798472113525944
$ws.Range("I6").Value = 0.2899264353016711
$ws.Range("J6").Value = 0.2899264353016712
$ws.Range("O6").Value = 0.3040809095127364
$ws.Range("P6").Value = 0.3040809095127364
$ws.Range("S6").Value = 0.08816109413831769
$ws.Range("T6").Value = 0.0881610941383177
$ws.Range("I7").Value = 0.2899264353016711
$ws.Range("J7").Value = 0.2899264353016712
$ws.Range("M7").Value = 2.009179666666667
$ws.Range("N7").Value = 6.027539
$ws.Range("O7").Value = 0.4205140165432039
$ws.Range("P7").Value = 0.4205140165432039
$ws.Range("Q7").Value = 1.125626165086111
$ws.Range("R7").Value = 10.130635485775
$ws.Range("S7").Value = 0.1219181298107591
$ws.Range("T7").Value = 0.1219181298107591
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.360498
$ws.Range("H8").Value = 1.081494
$ws.Range("I8").Value = 0.1865585983549632
$ws.Range("J8").Value = 0.1865585983549632
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.315861666666667
$ws.Range("N8").Value = 3.947585
$ws.Range("O8").Value = 0.2754050739440597
$ws.Range("P8").Value = 0.2754050739440597
$ws.Range("Q8").Value = 0.47436549911
$ws.Range("R8").Value = 4.26928949199
$ws.Range("S8").Value = 0.05137918457484877
$ws.Range("T8").Value = 0.05137918457484877
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.360498
$ws.Range("H9").Value = 1.081494
$ws.Range("I9").Value = 0.1865585983549632
$ws.Range("J9").Value = 0.1865585983549632
$ws.Range("O9").Value = 0.3040809095127364
$ws.Range("P9").Value = 0.3040809095127364
$ws.Range("Q9").Value = 0.523757570422
$ws.Range("R9").Value = 4.713818133797999
$ws.Range("S9").Value = 0.0567289082651985
$ws.Range("T9").Value = 0.0567289082651985
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.360498
$ws.Range("H10").Value = 1.081494
$ws.Range("I10").Value = 0.1865585983549632
$ws.Range("J10").Value = 0.1865585983549632
$ws.Range("M10").Value = 2.009179666666667
$ws.Range("N10").Value = 6.027539
$ws.Range("O10").Value = 0.4205140165432039
$ws.Range("P10").Value = 0.4205140165432039
$ws.Range("Q10").Value = 0.724305251474
$ws.Range("R10").Value = 6.518747263266
$ws.Range("S10").Value = 0.07845050551491592
$ws.Range("T10").Value = 0.07845050551491592
